$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '62.137.75'
$ws.Range('E2').Value = '  -6.51%  '
Set-TextValue $ws.Range('D3') '3.021.41'
$ws.Range('E3').Value = '  -7.14%  '
Set-TextValue $ws.Range('D4') '0.999'
$ws.Range('E4').Value = '  -0.13%  '
Set-TextValue $ws.Range('D5') '536.55'
$ws.Range('E5').Value = '  -8.02%  '
Set-TextValue $ws.Range('D6') '131.62'
$ws.Range('E6').Value = '  -14.59%  '
$ws.Range('E7').Value = '  +0.09%  '
Set-TextValue $ws.Range('D8') '3.017.75'
$ws.Range('E8').Value = '  -7.03%  '
Set-TextValue $ws.Range('D9') '0.482'
$ws.Range('E9').Value = '  -6.73%  '
$ws.Range('E10').Value = '  -8.28%  '
$ws.Range('E11').Value = '  -13.83%  '
Set-TextValue $ws.Range('D12') '0.456'
$ws.Range('E12').Value = '  -7.25%  '
Set-TextValue $ws.Range('D13') '34.09'
$ws.Range('E13').Value = '  -10.55%  '
$ws.Range('E14').Value = '  -10.93%  '
Set-TextValue $ws.Range('D15') '3.501.95'
$ws.Range('E15').Value = '  -7.33%  '
Set-TextValue $ws.Range('D16') '62.037.11'
$ws.Range('E16').Value = '  -6.79%  '
$ws.Range('E17').Value = '  -4.50%  '
Set-TextValue $ws.Range('D18') '3.040.20'
$ws.Range('E18').Value = '  -6.66%  '
Set-TextValue $ws.Range('D19') '6.54'
$ws.Range('E19').Value = '  -8.61%  '
Set-TextValue $ws.Range('D20') '473.54'
$ws.Range('E20').Value = '  -15.07%  '
Set-TextValue $ws.Range('D21') '13.20'
$ws.Range('E21').Value = '  -9.30%  '
Set-TextValue $ws.Range('D22') '0.696'
$ws.Range('E22').Value = '  -7.07%  '
Set-TextValue $ws.Range('D23') '7.05'
$ws.Range('E23').Value = '  -10.43%  '
Set-TextValue $ws.Range('D24') '77.03'
$ws.Range('E24').Value = '  -6.01%  '
$ws.Range('E25').Value = '  -12.51%  '
$ws.Range('E26').Value = '  -0.67%  '
$ws.Range('E27').Value = '  -10.60%  '
Set-TextValue $ws.Range('D28') '8.19'
$ws.Range('E28').Value = '  -12.00%  '
Set-TextValue $ws.Range('D29') '0.997'
$ws.Range('E29').Value = '  -0.41%  '
$ws.Range('E30').Value = '  -16.43%  '
Set-TextValue $ws.Range('D31') '25.83'
$ws.Range('E31').Value = '  -7.38%  '
$ws.Range('E32').Value = '  -8.22%  '
Set-TextValue $ws.Range('D33') '58.90'
$ws.Range('E34').Value = '  -13.35%  '
Set-TextValue $ws.Range('D35') '485.86'
$ws.Range('E35').Value = '  -14.25%  '
Set-TextValue $ws.Range('D36') '5.85'
$ws.Range('E36').Value = '  -9.17%  '
$ws.Range('E37').Value = '  -13.79%  '
Set-TextValue $ws.Range('D38') '3.098.73'
$ws.Range('E38').Value = '  -3.48%  '
Set-TextValue $ws.Range('D39') '0.0388'
$ws.Range('E39').Value = '  -14.74%  '
Set-TextValue $ws.Range('D40') '0.0779'
$ws.Range('E40').Value = '  -10.33%  '
$ws.Range('E41').Value = '  -12.56%  '
Set-TextValue $ws.Range('D42') '7.91'
$ws.Range('E42').Value = '  -8.83%  '
Set-TextValue $ws.Range('D43') '2.51'
$ws.Range('E43').Value = '  -17.41%  '
$ws.Range('E45').Value = '  -12.61%  '
$ws.Range('B46').Value = 'Fetch.AI'
$ws.Range('C46').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range('D46') '2.00'
$ws.Range('E46').Value = '  -13.81%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D47') '24.35'
$ws.Range('E47').Value = '  -8.28%  '
Set-TextValue $ws.Range('D48') '116.76'
$ws.Range('E48').Value = '  -7.81%  '
$ws.Range('E49').Value = '  -6.53%  '
Set-TextValue $ws.Range('D50') '0.0₃0492'
$ws.Range('E50').Value = '  -12.52%  '
Set-TextValue $ws.Range('D51') '2.27'
$ws.Range('E51').Value = '  +18.96%  '
